$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-CellText {
    param($cellRef, $val)
    $rng = $ws.Range($cellRef)
    $origStyle = $rng.Style
    $rng.NumberFormat = "@"
    $rng.Value = $val
    $rng.Style = $origStyle
}

$ws.Range("D2").Value = "33.998.06"
$ws.Range("E2").Value = "  -0.61%  "
$ws.Range("D3").Value = "1.786.26"
$ws.Range("E3").Value = "  -1.61%  "
$ws.Range("E4").Value = "  +0.20%  "
Set-CellText "D5" "227.22"
$ws.Range("E5").Value = "  +0.99%  "
Set-CellText "D6" "0.549"
$ws.Range("E6").Value = "  -1.44%  "
$ws.Range("E7").Value = "  +0.26%  "
Set-CellText "D8" "32.85"
$ws.Range("E8").Value = "  +2.29%  "
$ws.Range("E9").Value = "  -1.91%  "
$ws.Range("E10").Value = "  -3.77%  "
Set-CellText "D11" "0.0935"
$ws.Range("E11").Value = "  +0.67%  "
$ws.Range("D12").Value = "2.040.95"
$ws.Range("E12").Value = "  -1.63%  "
Set-CellText "D13" "11.27"
$ws.Range("E13").Value = "  +1.40%  "
$ws.Range("D14").Value = "1.782.01"
$ws.Range("E14").Value = "  -1.87%  "
$ws.Range("B15").Value = "Polygon"
$ws.Range("C15").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
Set-CellText "D15" "0.622"
$ws.Range("E15").Value = "  -3.31%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "33.959.96"
$ws.Range("E16").Value = "  -0.68%  "
Set-CellText "D17" "4.13"
$ws.Range("E17").Value = "  -4.66%  "
Set-CellText "D18" "67.95"
$ws.Range("E18").Value = "  -2.51%  "
Set-CellText "D19" "245.04"
$ws.Range("E19").Value = "  -2.05%  "
$ws.Range("D20").Value = "0.0₃0786"
$ws.Range("E20").Value = "  -2.55%  "
$ws.Range("E21").Value = "  +0.26%  "
Set-CellText "D22" "10.75"
$ws.Range("E22").Value = "  -3.21%  "
$ws.Range("E23").Value = "  -4.45%  "
Set-CellText "D24" "2.08"
$ws.Range("E24").Value = "  -3.81%  "
Set-CellText "D25" "160.04"
$ws.Range("E25").Value = "  -0.42%  "
Set-CellText "D26" "16.33"
$ws.Range("E26").Value = "  -2.29%  "
Set-CellText "D27" "7.08"
$ws.Range("E27").Value = "  -2.36%  "
$ws.Range("E28").Value = "  -2.31%  "
$ws.Range("E29").Value = "  +0.20%  "
$ws.Range("E30").Value = "  +1.50%  "
$ws.Range("E31").Value = "  -3.95%  "
Set-CellText "D32" "3.65"
$ws.Range("E32").Value = "  -3.71%  "
$ws.Range("E33").Value = "  -2.80%  "
$ws.Range("E34").Value = "  -4.65%  "
$ws.Range("D35").Value = "1.394.06"
$ws.Range("E35").Value = "  -2.72%  "
$ws.Range("E37").Value = "  -2.07%  "
$ws.Range("E38").Value = "  -2.15%  "
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText "D39" "2.21"
$ws.Range("E39").Value = "  +1.65%  "
$ws.Range("B40").Value = "HuobiToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"
Set-CellText "D40" "2.36"
$ws.Range("E40").Value = "  +0.93%  "
Set-CellText "D41" "0.916"
$ws.Range("E41").Value = "  -4.95%  "
$ws.Range("B42").Value = "Aave"
$ws.Range("C42").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-CellText "D42" "78.18"
$ws.Range("E42").Value = "  -4.47%  "
$ws.Range("B43").Value = "MXToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
Set-CellText "D43" "2.68"
$ws.Range("E43").Value = "  -2.73%  "
Set-CellText "D44" "13.04"
$ws.Range("E44").Value = "  +9.74%  "
$ws.Range("E45").Value = "  +11.06%  "
$ws.Range("E46").Value = "  +3.66%  "
Set-CellText "D47" "108.18"
$ws.Range("E47").Value = "  +1.69%  "
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("D49").Value = "1.940.30"
$ws.Range("E49").Value = "  -1.37%  "
Set-CellText "D50" "5.80"
$ws.Range("E50").Value = "  -4.60%  "
$ws.Range("E51").Value = "  +0.36%  "
